# Update vm_pu values for the Case with 380 kV (B2 slack voltage changed 1.05 -> 1.02,
# and recomputed per-unit voltages for the other buses), rows 2-25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.038618551531024
$ws.Cells.Item(2, 4).Value = 1.048378503482646
$ws.Cells.Item(2, 5).Value = 1.046995809523267
$ws.Cells.Item(2, 6).Value = 1.057323428498605
$ws.Cells.Item(2, 9).Value = 1.038661525066316
$ws.Cells.Item(2, 10).Value = 1.043715010298305
$ws.Cells.Item(2, 11).Value = 1.051138514899363
$ws.Cells.Item(2, 12).Value = 1.049759684168515
$ws.Cells.Item(2, 13).Value = 1.060058712381253
$ws.Cells.Item(2, 14).Value = 1.018381805078532
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.040076674357336
$ws.Cells.Item(3, 4).Value = 1.049157445373843
$ws.Cells.Item(3, 5).Value = 1.048225019816192
$ws.Cells.Item(3, 6).Value = 1.058483315928977
$ws.Cells.Item(3, 9).Value = 1.038855912981054
$ws.Cells.Item(3, 10).Value = 1.044815452581159
$ws.Cells.Item(3, 11).Value = 1.051729575819871
$ws.Cells.Item(3, 12).Value = 1.050799566311564
$ws.Cells.Item(3, 13).Value = 1.061031534212488
$ws.Cells.Item(3, 14).Value = 1.018755303345353
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.041019219042549
$ws.Cells.Item(4, 4).Value = 1.049659330631839
$ws.Cells.Item(4, 5).Value = 1.049019522073246
$ws.Cells.Item(4, 6).Value = 1.059232202927402
$ws.Cells.Item(4, 9).Value = 1.038978964061615
$ws.Cells.Item(4, 10).Value = 1.04552616350787
$ws.Cells.Item(4, 11).Value = 1.052109236075215
$ws.Cells.Item(4, 12).Value = 1.051471001559524
$ws.Cells.Item(4, 13).Value = 1.061658801469385
$ws.Cells.Item(4, 14).Value = 1.018996293964422
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.041415241111348
$ws.Cells.Item(5, 4).Value = 1.04986981087878
$ws.Cells.Item(5, 5).Value = 1.049353323950254
$ws.Cells.Item(5, 6).Value = 1.059546645587892
$ws.Cells.Item(5, 9).Value = 1.039030040994464
$ws.Cells.Item(5, 10).Value = 1.045824627307159
$ws.Cells.Item(5, 11).Value = 1.052268176577252
$ws.Cells.Item(5, 12).Value = 1.051752931895971
$ws.Cells.Item(5, 13).Value = 1.061921977050529
$ws.Cells.Item(5, 14).Value = 1.019097443023389
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.041481722005378
$ws.Cells.Item(6, 4).Value = 1.049905121410827
$ws.Cells.Item(6, 5).Value = 1.049409358726041
$ws.Cells.Item(6, 6).Value = 1.059599419133622
$ws.Cells.Item(6, 9).Value = 1.039038578701236
$ws.Cells.Item(6, 10).Value = 1.045874722092558
$ws.Cells.Item(6, 11).Value = 1.052294824193153
$ws.Cells.Item(6, 12).Value = 1.051800249334324
$ws.Cells.Item(6, 13).Value = 1.061966134518902
$ws.Cells.Item(6, 14).Value = 1.019114416849842
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.041024511580841
$ws.Cells.Item(7, 4).Value = 1.049662145091698
$ws.Cells.Item(7, 5).Value = 1.049023983158766
$ws.Cells.Item(7, 6).Value = 1.059236406052646
$ws.Cells.Item(7, 9).Value = 1.03897964912287
$ws.Cells.Item(7, 10).Value = 1.045530152844804
$ws.Cells.Item(7, 11).Value = 1.052111362472737
$ws.Cells.Item(7, 12).Value = 1.0514747700626
$ws.Cells.Item(7, 13).Value = 1.061662320104574
$ws.Cells.Item(7, 14).Value = 1.018997646163605
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.039111532826217
$ws.Cells.Item(8, 4).Value = 1.048642194442462
$ws.Cells.Item(8, 5).Value = 1.047411411052578
$ws.Cells.Item(8, 6).Value = 1.057715758591147
$ws.Cells.Item(8, 9).Value = 1.038727784855415
$ws.Cells.Item(8, 10).Value = 1.044087191518014
$ws.Cells.Item(8, 11).Value = 1.051338845956662
$ws.Cells.Item(8, 12).Value = 1.050111416912977
$ws.Cells.Item(8, 13).Value = 1.060387941822254
$ws.Cells.Item(8, 14).Value = 1.018508173720581
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.03573302413273
$ws.Cells.Item(9, 4).Value = 1.046828477299191
$ws.Cells.Item(9, 5).Value = 1.044562959859227
$ws.Cells.Item(9, 6).Value = 1.055023519883015
$ws.Cells.Item(9, 9).Value = 1.038263048910421
$ws.Cells.Item(9, 10).Value = 1.041534005952084
$ws.Cells.Item(9, 11).Value = 1.04995613115067
$ws.Cells.Item(9, 12).Value = 1.047697854874158
$ws.Cells.Item(9, 13).Value = 1.058125264235886
$ws.Cells.Item(9, 14).Value = 1.017640337156731
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.033475216925498
$ws.Cells.Item(10, 4).Value = 1.045608239752404
$ws.Cells.Item(10, 5).Value = 1.042659129073583
$ws.Cells.Item(10, 6).Value = 1.053220006089551
$ws.Cells.Item(10, 9).Value = 1.037939151716248
$ws.Cells.Item(10, 10).Value = 1.039824575445379
$ws.Cells.Item(10, 11).Value = 1.049019856844616
$ws.Cells.Item(10, 12).Value = 1.046081102873132
$ws.Cells.Item(10, 13).Value = 1.056605187145505
$ws.Cells.Item(10, 14).Value = 1.017058120219122
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.032496183379533
$ws.Cells.Item(11, 4).Value = 1.045077220819457
$ws.Cells.Item(11, 5).Value = 1.041833543822333
$ws.Cells.Item(11, 6).Value = 1.052436963838557
$ws.Cells.Item(11, 9).Value = 1.037795559672326
$ws.Cells.Item(11, 10).Value = 1.03908258469198
$ws.Cells.Item(11, 11).Value = 1.04861099660739
$ws.Cells.Item(11, 12).Value = 1.045379155329802
$ws.Cells.Item(11, 13).Value = 1.055944183612945
$ws.Cells.Item(11, 14).Value = 1.016805127796378
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.032132311093363
$ws.Cells.Item(12, 4).Value = 1.044879577724631
$ws.Cells.Item(12, 5).Value = 1.041526697720705
$ws.Cells.Item(12, 6).Value = 1.052145786619557
$ws.Cells.Item(12, 9).Value = 1.037741720646415
$ws.Cells.Item(12, 10).Value = 1.038806701572593
$ws.Cells.Item(12, 11).Value = 1.048458608550889
$ws.Cells.Item(12, 12).Value = 1.045118133634237
$ws.Cells.Item(12, 13).Value = 1.055698233518394
$ws.Cells.Item(12, 14).Value = 1.016711020030701
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.03221037280624
$ws.Cells.Item(13, 4).Value = 1.044921990915344
$ws.Cells.Item(13, 5).Value = 1.04159252578664
$ws.Cells.Item(13, 6).Value = 1.052208259681801
$ws.Cells.Item(13, 9).Value = 1.03775329204803
$ws.Cells.Item(13, 10).Value = 1.038865891948482
$ws.Cells.Item(13, 11).Value = 1.048491319819806
$ws.Cells.Item(13, 12).Value = 1.045174136734595
$ws.Cells.Item(13, 13).Value = 1.055751009886185
$ws.Cells.Item(13, 14).Value = 1.016731212606093
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.032466110022927
$ws.Cells.Item(14, 4).Value = 1.045060891726477
$ws.Cells.Item(14, 5).Value = 1.04180818366894
$ws.Cells.Item(14, 6).Value = 1.052412901597292
$ws.Cells.Item(14, 9).Value = 1.037791119575793
$ws.Cells.Item(14, 10).Value = 1.039059785739773
$ws.Cells.Item(14, 11).Value = 1.048598410761162
$ws.Cells.Item(14, 12).Value = 1.045357585076621
$ws.Cells.Item(14, 13).Value = 1.055923861968944
$ws.Cells.Item(14, 14).Value = 1.016797351584765
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.032623649375819
$ws.Cells.Item(15, 4).Value = 1.045146420190351
$ws.Cells.Item(15, 5).Value = 1.041941032692871
$ws.Cells.Item(15, 6).Value = 1.05253894566176
$ws.Cells.Item(15, 9).Value = 1.037814359759704
$ws.Cells.Item(15, 10).Value = 1.039179213554627
$ws.Cells.Item(15, 11).Value = 1.048664324192881
$ws.Cells.Item(15, 12).Value = 1.045470575482083
$ws.Cells.Item(15, 13).Value = 1.056030305565801
$ws.Cells.Item(15, 14).Value = 1.01683808403999
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.033540162307784
$ws.Cells.Item(16, 4).Value = 1.045643425839299
$ws.Cells.Item(16, 5).Value = 1.042713894479683
$ws.Cells.Item(16, 6).Value = 1.053271929284805
$ws.Cells.Item(16, 9).Value = 1.037948610993931
$ws.Cells.Item(16, 10).Value = 1.039873780671599
$ws.Cells.Item(16, 11).Value = 1.049046918816463
$ws.Cells.Item(16, 12).Value = 1.046127648785396
$ws.Cells.Item(16, 13).Value = 1.056648996500089
$ws.Cells.Item(16, 14).Value = 1.017074891633469
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.034114690029321
$ws.Cells.Item(17, 4).Value = 1.045954474579301
$ws.Cells.Item(17, 5).Value = 1.043198362017261
$ws.Cells.Item(17, 6).Value = 1.053731143394766
$ws.Cells.Item(17, 9).Value = 1.038031928120356
$ws.Cells.Item(17, 10).Value = 1.040308979985812
$ws.Cells.Item(17, 11).Value = 1.049285986514601
$ws.Cells.Item(17, 12).Value = 1.046539306365182
$ws.Cells.Item(17, 13).Value = 1.057036333015341
$ws.Cells.Item(17, 14).Value = 1.017223195734891
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.034449668877211
$ws.Cells.Item(18, 4).Value = 1.046135648634924
$ws.Cells.Item(18, 5).Value = 1.043480826919968
$ws.Cells.Item(18, 6).Value = 1.053998791988795
$ws.Cells.Item(18, 9).Value = 1.038080203123669
$ws.Cells.Item(18, 10).Value = 1.040562651421332
$ws.Cells.Item(18, 11).Value = 1.049425098132932
$ws.Cells.Item(18, 12).Value = 1.046779237845582
$ws.Cells.Item(18, 13).Value = 1.057261989904697
$ws.Cells.Item(18, 14).Value = 1.017309613376207
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.034563865585155
$ws.Cells.Item(19, 4).Value = 1.046197381018224
$ws.Cells.Item(19, 5).Value = 1.043577120476283
$ws.Cells.Item(19, 6).Value = 1.054090018880306
$ws.Cells.Item(19, 9).Value = 1.038096608967281
$ws.Cells.Item(19, 10).Value = 1.040649117634693
$ws.Cells.Item(19, 11).Value = 1.049472475234275
$ws.Cells.Item(19, 12).Value = 1.046861017633431
$ws.Cells.Item(19, 13).Value = 1.057338887382365
$ws.Cells.Item(19, 14).Value = 1.017339065082242
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.03405306249705
$ws.Cells.Item(20, 4).Value = 1.045921128420352
$ws.Cells.Item(20, 5).Value = 1.04314639535161
$ws.Cells.Item(20, 6).Value = 1.053681895113839
$ws.Cells.Item(20, 9).Value = 1.038023022338292
$ws.Cells.Item(20, 10).Value = 1.040262305143725
$ws.Cells.Item(20, 11).Value = 1.049260371216835
$ws.Cells.Item(20, 12).Value = 1.046495158198417
$ws.Cells.Item(20, 13).Value = 1.056994803418575
$ws.Cells.Item(20, 14).Value = 1.01720729297733
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.032390807796464
$ws.Cells.Item(21, 4).Value = 1.045019999924234
$ws.Cells.Item(21, 5).Value = 1.041744682996438
$ws.Cells.Item(21, 6).Value = 1.052352648529444
$ws.Cells.Item(21, 9).Value = 1.037779994187689
$ws.Cells.Item(21, 10).Value = 1.039002696473821
$ws.Cells.Item(21, 11).Value = 1.048566889489663
$ws.Cells.Item(21, 12).Value = 1.045303572072066
$ws.Cells.Item(21, 13).Value = 1.055872973068386
$ws.Cells.Item(21, 14).Value = 1.016777879046834
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.03134443114157
$ws.Cells.Item(22, 4).Value = 1.044451115449395
$ws.Cells.Item(22, 5).Value = 1.04086228668622
$ws.Cells.Item(22, 6).Value = 1.051515042516546
$ws.Cells.Item(22, 9).Value = 1.037624284949395
$ws.Cells.Item(22, 10).Value = 1.038209138637362
$ws.Cells.Item(22, 11).Value = 1.048127865656787
$ws.Cells.Item(22, 12).Value = 1.04455271126727
$ws.Cells.Item(22, 13).Value = 1.055165179297903
$ws.Cells.Item(22, 14).Value = 1.016507106923383
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.031899256177548
$ws.Cells.Item(23, 4).Value = 1.044752911167421
$ws.Cells.Item(23, 5).Value = 1.041330165973226
$ws.Cells.Item(23, 6).Value = 1.051959250455753
$ws.Cells.Item(23, 9).Value = 1.03770710512342
$ws.Cells.Item(23, 10).Value = 1.038629971191248
$ws.Cells.Item(23, 11).Value = 1.048360885700667
$ws.Cells.Item(23, 12).Value = 1.044950915857172
$ws.Cells.Item(23, 13).Value = 1.055540627927925
$ws.Cells.Item(23, 14).Value = 1.016650723064882
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.034080909738855
$ws.Cells.Item(24, 4).Value = 1.045936196903871
$ws.Cells.Item(24, 5).Value = 1.043169877212285
$ws.Cells.Item(24, 6).Value = 1.053704148920543
$ws.Cells.Item(24, 9).Value = 1.038027047474467
$ws.Cells.Item(24, 10).Value = 1.040283396030528
$ws.Cells.Item(24, 11).Value = 1.049271946695149
$ws.Cells.Item(24, 12).Value = 1.046515107415565
$ws.Cells.Item(24, 13).Value = 1.057013569690259
$ws.Cells.Item(24, 14).Value = 1.017214479013665
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.036607386761951
$ws.Cells.Item(25, 4).Value = 1.047299318465106
$ws.Cells.Item(25, 5).Value = 1.045300193046902
$ws.Cells.Item(25, 6).Value = 1.055721046531371
$ws.Cells.Item(25, 9).Value = 1.038385674070529
$ws.Cells.Item(25, 10).Value = 1.042195335400377
$ws.Cells.Item(25, 11).Value = 1.050316141844271
$ws.Cells.Item(25, 12).Value = 1.048323162213575
$ws.Cells.Item(25, 13).Value = 1.058712258145474
$ws.Cells.Item(25, 14).Value = 1.017865332767142
